# TDoA estimated uncertainty monte carlo runs
# Update the Monte Carlo simulation result table (rows 18-27, columns B:H)
# with the refreshed run values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 18
$ws.Range("B18").Value = 0.46141975308642003
$ws.Range("C18").Value = 7.3443976533006996
$ws.Range("D18").Value = 11.878178357748199
$ws.Range("E18").Value = 5.3650685428053499
$ws.Range("F18").Value = 10.200537333616699
$ws.Range("G18").Value = 50
$ws.Range("H18").Value = 45

# Row 19
$ws.Range("B19").Value = 0.54783950617284005
$ws.Range("C19").Value = 5.7331995097676698
$ws.Range("D19").Value = 9.0782550527084105
$ws.Range("E19").Value = 4.1838353121013698
$ws.Range("F19").Value = 7.7017234041168203
$ws.Range("G19").Value = 40
$ws.Range("H19").Value = 55

# Row 20 (B/G/H unchanged)
$ws.Range("C20").Value = 123.29050740973
$ws.Range("D20").Value = 100.175116397633
$ws.Range("E20").Value = 42.7775199530858
$ws.Range("F20").Value = 7.1602831439958301

# Row 21
$ws.Range("B21").Value = 0.407407407407407
$ws.Range("C21").Value = 8.1225063682991703
$ws.Range("D21").Value = 12.4190159570666
$ws.Range("E21").Value = 6.0103433874487404
$ws.Range("F21").Value = 13.2229813561994
$ws.Range("G21").Value = 55
$ws.Range("H21").Value = 45

# Row 22 (G/H unchanged)
$ws.Range("B22").Value = 0.101851851851852
$ws.Range("C22").Value = 11.9640504941136
$ws.Range("D22").Value = 15.6073273723034
$ws.Range("E22").Value = 9.4402349375060606
$ws.Range("F22").Value = 17.670167685787899

# Row 23
$ws.Range("B23").Value = 0.58950617283950602
$ws.Range("C23").Value = 4.9482003274105901
$ws.Range("D23").Value = 8.2595696927181397
$ws.Range("E23").Value = 3.6907206883117598
$ws.Range("F23").Value = 6.7197458555538798
$ws.Range("G23").Value = 40
$ws.Range("H23").Value = 60

# Row 24 (G/H unchanged)
$ws.Range("B24").Value = 0.19290123456790101
$ws.Range("C24").Value = 11.1082214821676
$ws.Range("D24").Value = 14.3663316054086
$ws.Range("E24").Value = 8.3116644001215398
$ws.Range("F24").Value = 16.375307989934502

# Row 25 (highlighted best-pick row)
$ws.Range("B25").Value = 0.60030864197530898
$ws.Range("C25").Value = 4.8018626749266398
$ws.Range("D25").Value = 8.3406348000448602
$ws.Range("E25").Value = 3.5278827670437698
$ws.Range("F25").Value = 6.1454064014607699
$ws.Range("G25").Value = 35
$ws.Range("H25").Value = 65

# Row 26
$ws.Range("B26").Value = 0.57407407407407396
$ws.Range("C26").Value = 5.3467429659569303
$ws.Range("D26").Value = 8.3074194251262803
$ws.Range("E26").Value = 3.9839773502531002
$ws.Range("F26").Value = 7.0424026094141796
$ws.Range("G26").Value = 40
$ws.Range("H26").Value = 60

# Row 27
$ws.Range("B27").Value = 0.45061728395061701
$ws.Range("C27").Value = 7.4607002139715402
$ws.Range("D27").Value = 11.3092209729366
$ws.Range("E27").Value = 5.52490595590643
$ws.Range("F27").Value = 10.874135367838701
$ws.Range("G27").Value = 50
$ws.Range("H27").Value = 50

# Update the window scroll position and selection to match where the author
# left off reviewing the refreshed table.
$ws.Activate()
$ws.Range("B27:H27").Select()
$excel.ActiveWindow.ScrollRow = 12
$excel.ActiveWindow.ScrollColumn = 1
